$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Cell="D2"; Value="307.49"},
    @{Cell="E2"; Value="0.66%"},
    @{Cell="D3"; Value="38.56"},
    @{Cell="E3"; Value="8.32%"},
    @{Cell="D4"; Value="5.101"},
    @{Cell="E4"; Value="1.15%"},
    @{Cell="E5"; Value="1.26%"},
    @{Cell="D6"; Value="1.965"},
    @{Cell="E6"; Value="4.78%"},
    @{Cell="D7"; Value="7.948"},
    @{Cell="E7"; Value="2.03%"},
    @{Cell="D8"; Value="0.9293"},
    @{Cell="E8"; Value="0.80%"},
    @{Cell="D9"; Value="0.1430"},
    @{Cell="E9"; Value="11.18%"},
    @{Cell="D10"; Value="0.1958"},
    @{Cell="E10"; Value="3.10%"},
    @{Cell="D11"; Value="0.09080"},
    @{Cell="E11"; Value="-0.47%"},
    @{Cell="D12"; Value="0.03509"},
    @{Cell="E12"; Value="2.70%"},
    @{Cell="D13"; Value="0.09827"},
    @{Cell="E13"; Value="-0.36%"},
    @{Cell="E14"; Value="0.41%"},
    @{Cell="D15"; Value="0.006168"},
    @{Cell="E15"; Value="-0.46%"},
    @{Cell="E16"; Value="-3.17%"},
    @{Cell="D17"; Value="4.195"},
    @{Cell="E18"; Value="1.99%"},
    @{Cell="D19"; Value="0.3463"},
    @{Cell="E19"; Value="1.24%"},
    @{Cell="E20"; Value="-0.47%"},
    @{Cell="D21"; Value="4.791"},
    @{Cell="E21"; Value="-7.64%"},
    @{Cell="E22"; Value="6.30%"},
    @{Cell="D23"; Value="0.04363"},
    @{Cell="E23"; Value="-1.35%"},
    @{Cell="D24"; Value="0.001223"},
    @{Cell="E24"; Value="-0.87%"},
    @{Cell="E25"; Value="-1.06%"},
    @{Cell="E27"; Value="4.06%"},
    @{Cell="D39"; Value="0.02086"},
    @{Cell="E39"; Value="7.73%"},
    @{Cell="D40"; Value="0.05122"},
    @{Cell="E40"; Value="-0.89%"},
    @{Cell="E41"; Value="-1.59%"},
    @{Cell="D42"; Value="0.01013"},
    @{Cell="E42"; Value="-0.11%"},
    @{Cell="D43"; Value="0.1355"},
    @{Cell="E43"; Value="0.23%"},
    @{Cell="D44"; Value="0.002133"},
    @{Cell="E44"; Value="-0.87%"},
    @{Cell="D45"; Value="0.009283"},
    @{Cell="E45"; Value="-3.57%"},
    @{Cell="D46"; Value="0.00006245"},
    @{Cell="E47"; Value="0.00%"},
    @{Cell="D48"; Value="0.003028"},
    @{Cell="E49"; Value="-3.58%"},
    @{Cell="D50"; Value="0.00002103"},
    @{Cell="E50"; Value="0.00%"},
    @{Cell="D51"; Value="0.0002003"},
    @{Cell="E51"; Value="0.00%"}
)

foreach ($change in $changes) {
    $ws.Range("ZZ1").NumberFormat = "@"
    $ws.Range("ZZ1").Value = $change.Value
    $ws.Range("ZZ1").Copy()
    $ws.Range($change.Cell).PasteSpecial(-4163)
}

$ws.Range("ZZ1").Clear()
$excel.CutCopyMode = $false
